$wb = $excel.ActiveWorkbook

# 1) verifySearchWebinarsHistorical
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws17 = $wb.Worksheets.Add($null, $lastSheet)
$ws17.Name = "verifySearchWebinarsHistorical"
$ws17.Range("A1").Value = "TypeHistoricalSearchBar"
$ws17.Range("A2").Value = "Wood"
$ws17.Columns.Item(1).ColumnWidth = 22.7109375
[void]$ws17.Range("A1:A2").Select()

# 2) verifyInvalidSearchHistorical
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws18 = $wb.Worksheets.Add($null, $lastSheet)
$ws18.Name = "verifyInvalidSearchHistorical"
$ws18.Range("A1").Value = "TypeHistoricalSearchBarInvalid"
$ws18.Range("A2").Value = "AutomatedTesting"
[void]$ws18.Range("A3").Select()

# 3) verifySortWebinarsHistorical
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws19 = $wb.Worksheets.Add($null, $lastSheet)
$ws19.Name = "verifySortWebinarsHistorical"
$ws19.Range("A1").Value = "TypeHistoricalSearchBarSort"
$ws19.Range("A2").Value = "Wood"
$ws19.Columns.Item(1).ColumnWidth = 22.7109375
[void]$ws19.Range("A1:A2").Select()

# 4) verifyClearSearchHistorical
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws20 = $wb.Worksheets.Add($null, $lastSheet)
$ws20.Name = "verifyClearSearchHistorical"
$ws20.Range("A1").Value = "TypeHistoricalClearSearch"
$ws20.Range("A2").Value = "Wood"
$ws20.Columns.Item(1).ColumnWidth = 26.42578125
[void]$ws20.Range("B5").Select()

# 5) verifyExternalVideo (becomes the active / selected sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws21 = $wb.Worksheets.Add($null, $lastSheet)
$ws21.Name = "verifyExternalVideo"
$ws21.Range("A1").Value = "TypeExternal"
$ws21.Range("A2").Value = "Wood"
$ws21.Activate()
[void]$ws21.Range("A1").Select()
